# "Actualizacion asociacion usuario entrevista"
# The "Reclutador" column (C) is repurposed into a "Usuario" column that
# stores the login of the user associated with the interview record,
# instead of a manually typed recruiter name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the column header in C1 from "Reclutador" to "Usuario".
$ws.Range("C1").Value = "Usuario"

# Replace each row's recruiter name with the constant "login" value.
$ws.Range("C2").Value = "login"
$ws.Range("C3").Value = "login"
$ws.Range("C4").Value = "login"

# Re-apply the (unchanged) font explicitly so the cell style collapses
# onto the shared "Normal text" style used elsewhere in the sheet.
$ws.Range("C2:C4").Font.Name = "Calibri"
$ws.Range("C2:C4").Font.Size = 14
